$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 25,20
$data[0,0] = "ECs"
$data[0,1] = "Gnai2"
$data[0,2] = "Egfr"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 151.7260716666667
$data[0,7] = 455.178215
$data[0,8] = 0.2700739458961593
$data[0,9] = 0.2783366498663096
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 1.211057333333333
$data[0,13] = 3.633172
$data[0,14] = 0.0147461456544675
$data[0,15] = 0.01598314554371009
$data[0,16] = 183.7489717497756
$data[0,17] = 1653.74074574798
$data[0,18] = 0.00398254974366154
$data[0,19] = 0.004448695184961902
$data[1,0] = "ECs"
$data[1,1] = "Gnai2"
$data[1,2] = "Egfr"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 151.7260716666667
$data[1,7] = 455.178215
$data[1,8] = 0.2700739458961593
$data[1,9] = 0.2783366498663096
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 61.06015933333333
$data[1,13] = 183.180478
$data[1,14] = 0.7434842087418319
$data[1,15] = 0.8058523627949308
$data[1,16] = 9264.418110987419
$data[1,17] = 83379.76299888677
$data[1,18] = 0.2007957139663903
$data[1,19] = 0.2242982469471909
$data[2,0] = "ECs"
$data[2,1] = "Gnai2"
$data[2,2] = "Egfr"
$data[2,3] = "M1"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 151.7260716666667
$data[2,7] = 455.178215
$data[2,8] = 0.2700739458961593
$data[2,9] = 0.2783366498663096
$data[2,10] = 1
$data[2,11] = 0.3333333333333333
$data[2,12] = 0.132278
$data[2,13] = 0.396834
$data[2,14] = 0.001610650958623747
$data[2,15] = 0.001745762539921768
$data[2,16] = 20.07002130792334
$data[2,17] = 180.63019177131
$data[2,18] = 0.0004349948598569468
$data[2,19] = 0.0004859096968239245
$data[3,0] = "ECs"
$data[3,1] = "Gnai2"
$data[3,2] = "Egfr"
$data[3,3] = "M2"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 151.7260716666667
$data[3,7] = 455.178215
$data[3,8] = 0.2700739458961593
$data[3,9] = 0.2783366498663096
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 0.655122
$data[3,13] = 1.965366
$data[3,14] = 0.007976934012575832
$data[3,15] = 0.008646089649666828
$data[3,16] = 99.39908752241001
$data[3,17] = 894.59178770169
$data[3,18] = 0.002154362044929638
$data[3,19] = 0.002406523627532039
$data[4,0] = "ECs"
$data[4,1] = "Gnai2"
$data[4,2] = "Egfr"
$data[4,3] = "sCs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 151.7260716666667
$data[4,7] = 455.178215
$data[4,8] = 0.2700739458961593
$data[4,9] = 0.2783366498663096
$data[4,10] = 2
$data[4,11] = 1
$data[4,12] = 19.068426
$data[4,13] = 38.136852
$data[4,14] = 0.232182060632501
$data[4,15] = 0.1677726394717705
$data[4,16] = 2893.17736984653
$data[4,17] = 17359.06421907918
$data[4,18] = 0.06270632528132086
$data[4,19] = 0.04669727440980076
$data[5,0] = "FAPs"
$data[5,1] = "Gnai2"
$data[5,2] = "Egfr"
$data[5,3] = "ECs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 82.248871
$data[5,7] = 246.746613
$data[5,8] = 0.146403824289839
$data[5,9] = 0.150882936320401
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 1.211057333333333
$data[5,13] = 3.633172
$data[5,14] = 0.0147461456544675
$data[5,15] = 0.01598314554371009
$data[5,16] = 99.60809838293733
$data[5,17] = 896.472885446436
$data[5,18] = 0.002158892117349032
$data[5,19] = 0.002411583931271311
$data[6,0] = "FAPs"
$data[6,1] = "Gnai2"
$data[6,2] = "Egfr"
$data[6,3] = "FAPs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 82.248871
$data[6,7] = 246.746613
$data[6,8] = 0.146403824289839
$data[6,9] = 0.150882936320401
$data[6,10] = 3
$data[6,11] = 1
$data[6,12] = 61.06015933333333
$data[6,13] = 183.180478
$data[6,14] = 0.7434842087418319
$data[6,15] = 0.8058523627949308
$data[6,16] = 5022.129168246779
$data[6,17] = 45199.16251422101
$data[6,18] = 0.1088489314589091
$data[6,19] = 0.1215893707392322
$data[7,0] = "FAPs"
$data[7,1] = "Gnai2"
$data[7,2] = "Egfr"
$data[7,3] = "M1"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 82.248871
$data[7,7] = 246.746613
$data[7,8] = 0.146403824289839
$data[7,9] = 0.150882936320401
$data[7,10] = 1
$data[7,11] = 0.3333333333333333
$data[7,12] = 0.132278
$data[7,13] = 0.396834
$data[7,14] = 0.001610650958623747
$data[7,15] = 0.001745762539921768
$data[7,16] = 10.879716158138
$data[7,17] = 97.917445423242
$data[7,18] = 0.0002358054599386117
$data[7,19] = 0.0002634057781415576
$data[8,0] = "FAPs"
$data[8,1] = "Gnai2"
$data[8,2] = "Egfr"
$data[8,3] = "M2"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 82.248871
$data[8,7] = 246.746613
$data[8,8] = 0.146403824289839
$data[8,9] = 0.150882936320401
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 0.655122
$data[8,13] = 1.965366
$data[8,14] = 0.007976934012575832
$data[8,15] = 0.008646089649666828
$data[8,16] = 53.88304486726199
$data[8,17] = 484.947403805358
$data[8,18] = 0.001167853645548792
$data[8,19] = 0.001304547394031158
$data[9,0] = "FAPs"
$data[9,1] = "Gnai2"
$data[9,2] = "Egfr"
$data[9,3] = "sCs"
$data[9,4] = 3
$data[9,5] = 1
$data[9,6] = 82.248871
$data[9,7] = 246.746613
$data[9,8] = 0.146403824289839
$data[9,9] = 0.150882936320401
$data[9,10] = 2
$data[9,11] = 1
$data[9,12] = 19.068426
$data[9,13] = 38.136852
$data[9,14] = 0.232182060632501
$data[9,15] = 0.1677726394717705
$data[9,16] = 1568.356510247046
$data[9,17] = 9410.139061482276
$data[9,18] = 0.03399234160809342
$data[9,19] = 0.02531402847772473
$data[10,0] = "M1"
$data[10,1] = "Gnai2"
$data[10,2] = "Egfr"
$data[10,3] = "ECs"
$data[10,4] = 3
$data[10,5] = 1
$data[10,6] = 123.444321
$data[10,7] = 370.332963
$data[10,8] = 0.2197321429647646
$data[10,9] = 0.2264546783208506
$data[10,10] = 3
$data[10,11] = 1
$data[10,12] = 1.211057333333333
$data[10,13] = 3.633172
$data[10,14] = 0.0147461456544675
$data[10,15] = 0.01598314554371009
$data[10,16] = 149.498150205404
$data[10,17] = 1345.483351848636
$data[10,18] = 0.003240202185126695
$data[10,19] = 0.003619458082656206
$data[11,0] = "M1"
$data[11,1] = "Gnai2"
$data[11,2] = "Egfr"
$data[11,3] = "FAPs"
$data[11,4] = 3
$data[11,5] = 1
$data[11,6] = 123.444321
$data[11,7] = 370.332963
$data[11,8] = 0.2197321429647646
$data[11,9] = 0.2264546783208506
$data[11,10] = 3
$data[11,11] = 1
$data[11,12] = 61.06015933333333
$data[11,13] = 183.180478
$data[11,14] = 0.7434842087418319
$data[11,15] = 0.8058523627949308
$data[11,16] = 7537.529909055146
$data[11,17] = 67837.76918149632
$data[11,18] = 0.1633673784473051
$data[11,19] = 0.1824890375908234
$data[12,0] = "M1"
$data[12,1] = "Gnai2"
$data[12,2] = "Egfr"
$data[12,3] = "M1"
$data[12,4] = 3
$data[12,5] = 1
$data[12,6] = 123.444321
$data[12,7] = 370.332963
$data[12,8] = 0.2197321429647646
$data[12,9] = 0.2264546783208506
$data[12,10] = 1
$data[12,11] = 0.3333333333333333
$data[12,12] = 0.132278
$data[12,13] = 0.396834
$data[12,14] = 0.001610650958623747
$data[12,15] = 0.001745762539921768
$data[12,16] = 16.328967893238
$data[12,17] = 146.960711039142
$data[12,18] = 0.0003539117867066482
$data[12,19] = 0.0003953360944025751
$data[13,0] = "M1"
$data[13,1] = "Gnai2"
$data[13,2] = "Egfr"
$data[13,3] = "M2"
$data[13,4] = 3
$data[13,5] = 1
$data[13,6] = 123.444321
$data[13,7] = 370.332963
$data[13,8] = 0.2197321429647646
$data[13,9] = 0.2264546783208506
$data[13,10] = 3
$data[13,11] = 1
$data[13,12] = 0.655122
$data[13,13] = 1.965366
$data[13,14] = 0.007976934012575832
$data[13,15] = 0.008646089649666828
$data[13,16] = 80.871090462162
$data[13,17] = 727.839814159458
$data[13,18] = 0.001752788804871806
$data[13,19] = 0.001957947450348537
$data[14,0] = "M1"
$data[14,1] = "Gnai2"
$data[14,2] = "Egfr"
$data[14,3] = "sCs"
$data[14,4] = 3
$data[14,5] = 1
$data[14,6] = 123.444321
$data[14,7] = 370.332963
$data[14,8] = 0.2197321429647646
$data[14,9] = 0.2264546783208506
$data[14,10] = 2
$data[14,11] = 1
$data[14,12] = 19.068426
$data[14,13] = 38.136852
$data[14,14] = 0.232182060632501
$data[14,15] = 0.1677726394717705
$data[14,16] = 2353.888900108746
$data[14,17] = 14123.33340065247
$data[14,18] = 0.05101786174075436
$data[14,19] = 0.03799289910261982
$data[15,0] = "M2"
$data[15,1] = "Gnai2"
$data[15,2] = "Egfr"
$data[15,3] = "ECs"
$data[15,4] = 3
$data[15,5] = 1
$data[15,6] = 154.3429766666667
$data[15,7] = 463.02893
$data[15,8] = 0.2747320633285943
$data[15,9] = 0.2831372788071194
$data[15,10] = 3
$data[15,11] = 1
$data[15,12] = 1.211057333333333
$data[15,13] = 3.633172
$data[15,14] = 0.0147461456544675
$data[15,15] = 0.01598314554371009
$data[15,16] = 186.9181937406622
$data[15,17] = 1682.26374366596
$data[15,18] = 0.00405123902179584
$data[15,19] = 0.004525424336024213
$data[16,0] = "M2"
$data[16,1] = "Gnai2"
$data[16,2] = "Egfr"
$data[16,3] = "FAPs"
$data[16,4] = 3
$data[16,5] = 1
$data[16,6] = 154.3429766666667
$data[16,7] = 463.02893
$data[16,8] = 0.2747320633285943
$data[16,9] = 0.2831372788071194
$data[16,10] = 3
$data[16,11] = 1
$data[16,12] = 61.06015933333333
$data[16,13] = 183.180478
$data[16,14] = 0.7434842087418319
$data[16,15] = 0.8058523627949308
$data[16,16] = 9424.206747247616
$data[16,17] = 84817.86072522854
$data[16,18] = 0.2042589507198708
$data[16,19] = 0.2281668451220443
$data[17,0] = "M2"
$data[17,1] = "Gnai2"
$data[17,2] = "Egfr"
$data[17,3] = "M1"
$data[17,4] = 3
$data[17,5] = 1
$data[17,6] = 154.3429766666667
$data[17,7] = 463.02893
$data[17,8] = 0.2747320633285943
$data[17,9] = 0.2831372788071194
$data[17,10] = 1
$data[17,11] = 0.3333333333333333
$data[17,12] = 0.132278
$data[17,13] = 0.396834
$data[17,14] = 0.001610650958623747
$data[17,15] = 0.001745762539921768
$data[17,16] = 20.41618026751333
$data[17,17] = 183.74562240762
$data[17,18] = 0.0004424974611648802
$data[17,19] = 0.0004942904549968546
$data[18,0] = "M2"
$data[18,1] = "Gnai2"
$data[18,2] = "Egfr"
$data[18,3] = "M2"
$data[18,4] = 3
$data[18,5] = 1
$data[18,6] = 154.3429766666667
$data[18,7] = 463.02893
$data[18,8] = 0.2747320633285943
$data[18,9] = 0.2831372788071194
$data[18,10] = 3
$data[18,11] = 1
$data[18,12] = 0.655122
$data[18,13] = 1.965366
$data[18,14] = 0.007976934012575832
$data[18,15] = 0.008646089649666828
$data[18,16] = 101.11347955982
$data[18,17] = 910.02131603838
$data[18,18] = 0.002191519540311002
$data[18,19] = 0.002448030295729066
$data[19,0] = "M2"
$data[19,1] = "Gnai2"
$data[19,2] = "Egfr"
$data[19,3] = "sCs"
$data[19,4] = 3
$data[19,5] = 1
$data[19,6] = 154.3429766666667
$data[19,7] = 463.02893
$data[19,8] = 0.2747320633285943
$data[19,9] = 0.2831372788071194
$data[19,10] = 2
$data[19,11] = 1
$data[19,12] = 19.068426
$data[19,13] = 38.136852
$data[19,14] = 0.232182060632501
$data[19,15] = 0.1677726394717705
$data[19,16] = 2943.07762918806
$data[19,17] = 17658.46577512836
$data[19,18] = 0.06378785658545179
$data[19,19] = 0.047502688598325
$data[20,0] = "sCs"
$data[20,1] = "Gnai2"
$data[20,2] = "Egfr"
$data[20,3] = "ECs"
$data[20,4] = 2
$data[20,5] = 1
$data[20,6] = 50.0323125
$data[20,7] = 100.064625
$data[20,8] = 0.08905802352064279
$data[20,9] = 0.06118845668531954
$data[20,10] = 3
$data[20,11] = 1
$data[20,12] = 1.211057333333333
$data[20,13] = 3.633172
$data[20,14] = 0.0147461456544675
$data[20,15] = 0.01598314554371009
$data[20,16] = 60.59199895675001
$data[20,17] = 363.5519937405001
$data[20,18] = 0.001313262586534391
$data[20,19] = 0.000977984008796463
$data[21,0] = "sCs"
$data[21,1] = "Gnai2"
$data[21,2] = "Egfr"
$data[21,3] = "FAPs"
$data[21,4] = 2
$data[21,5] = 1
$data[21,6] = 50.0323125
$data[21,7] = 100.064625
$data[21,8] = 0.08905802352064279
$data[21,9] = 0.06118845668531954
$data[21,10] = 3
$data[21,11] = 1
$data[21,12] = 61.06015933333333
$data[21,13] = 183.180478
$data[21,14] = 0.7434842087418319
$data[21,15] = 0.8058523627949308
$data[21,16] = 3054.980973065125
$data[21,17] = 18329.88583839075
$data[21,18] = 0.06621323414935655
$data[21,19] = 0.04930886239564003
$data[22,0] = "sCs"
$data[22,1] = "Gnai2"
$data[22,2] = "Egfr"
$data[22,3] = "M1"
$data[22,4] = 2
$data[22,5] = 1
$data[22,6] = 50.0323125
$data[22,7] = 100.064625
$data[22,8] = 0.08905802352064279
$data[22,9] = 0.06118845668531954
$data[22,10] = 1
$data[22,11] = 0.3333333333333333
$data[22,12] = 0.132278
$data[22,13] = 0.396834
$data[22,14] = 0.001610650958623747
$data[22,15] = 0.001745762539921768
$data[22,16] = 6.618174232875001
$data[22,17] = 39.70904539725
$data[22,18] = 0.0001434413909566595
$data[22,19] = 0.0001068205155568565
$data[23,0] = "sCs"
$data[23,1] = "Gnai2"
$data[23,2] = "Egfr"
$data[23,3] = "M2"
$data[23,4] = 2
$data[23,5] = 1
$data[23,6] = 50.0323125
$data[23,7] = 100.064625
$data[23,8] = 0.08905802352064279
$data[23,9] = 0.06118845668531954
$data[23,10] = 3
$data[23,11] = 1
$data[23,12] = 0.655122
$data[23,13] = 1.965366
$data[23,14] = 0.007976934012575832
$data[23,15] = 0.008646089649666828
$data[23,16] = 32.777268629625
$data[23,17] = 196.66361177775
$data[23,18] = 0.000710409976914594
$data[23,19] = 0.0005290408820260283
$data[24,0] = "sCs"
$data[24,1] = "Gnai2"
$data[24,2] = "Egfr"
$data[24,3] = "sCs"
$data[24,4] = 2
$data[24,5] = 1
$data[24,6] = 50.0323125
$data[24,7] = 100.064625
$data[24,8] = 0.08905802352064279
$data[24,9] = 0.06118845668531954
$data[24,10] = 2
$data[24,11] = 1
$data[24,12] = 19.068426
$data[24,13] = 38.136852
$data[24,14] = 0.232182060632501
$data[24,15] = 0.1677726394717705
$data[24,16] = 954.037448515125
$data[24,17] = 3816.1497940605
$data[24,18] = 0.02067767541688059
$data[24,19] = 0.01026574888330016

$ws.Range("A2:T26").Value = $data
